$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab from "1" to "Dusheti"
$ws.Name = "Dusheti"

# Normalize the ellipsis placeholder character ("…") used for confidential /
# unavailable data to three literal dots ("...") everywhere it appears.
$ws.Cells.Replace("…", "...") | Out-Null

# The Urban and Rural breakdown rows are no longer published for any year -
# replace every year's value in those two rows with the "..." placeholder,
# matching the Total row's existing treatment of withheld data.
$ws.Range("B6:O7").Value = "..."

# Row 8 was a blank spacer row between the data table and the footnote -
# remove it so the footnote moves up to row 8.
$ws.Rows.Item(8).Delete()
